# Daily attendance processing - 2025-12-23 07:56:54
#
# In the "Recorded By" column (G), whenever the comma-separated list of
# recorders includes the literal token "System" (case-sensitive), move it to
# the front of the list while preserving the relative order of the other
# tokens. Rows whose "Recorded By" value doesn't contain "System" (or is just
# "System" alone) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $value = $cell.Value2

    if ($value -eq $null) {
        continue
    }

    $text = [string]$value
    if ($text -eq "") {
        continue
    }

    $parts = @($text -split ", ")

    $systemFound = $false
    $rest = @()
    foreach ($part in $parts) {
        if ($part.Equals("System")) {
            $systemFound = $true
        } else {
            $rest += $part
        }
    }

    if ($systemFound -and $rest.Count -gt 0) {
        $newParts = @("System") + $rest
        $newText = $newParts -join ", "
        if (-not $newText.Equals($text)) {
            $cell.Value = $newText
        }
    }
}
